$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs","Bmp6","Bmpr1b","FAPs",2,0.6666666666666666,4.145393333333334,12.43618,0.1621900462138432,0.1621900462138432,3,1,1.292736333333333,3.878209,0.9024289107934966,0.9024289107934966,5.358900577957778,48.23010520162,0.1463649867463054,0.1463649867463054),
    @("ECs","Bmp6","Bmpr1b","sCs",2,0.6666666666666666,4.145393333333334,12.43618,0.1621900462138432,0.1621900462138432,2,0.6666666666666666,0.1397713333333333,0.419314,0.09757108920650338,0.09757108920650338,0.5794071533911112,5.21466438052,0.0158250594675378,0.0158250594675378),
    @("FAPs","Bmp6","Bmpr1b","FAPs",3,1,6.772999666666666,20.318999,0.2649961151116367,0.2649961151116367,3,1,1.292736333333333,3.878209,0.9024289107934966,0.9024289107934966,8.755702754754555,78.80132479279099,0.2391401555247023,0.2391401555247023),
    @("FAPs","Bmp6","Bmpr1b","sCs",3,1,6.772999666666666,20.318999,0.2649961151116367,0.2649961151116367,2,0.6666666666666666,0.1397713333333333,0.419314,0.09757108920650338,0.09757108920650338,0.9466711940762221,8.520040746686,0.02585595958693434,0.02585595958693434),
    @("sCs","Bmp6","Bmpr1b","FAPs",3,1,14.640471,43.921413,0.5728138386745201,0.5728138386745202,3,1,1.292736333333333,3.878209,0.9024289107934966,0.9024289107934966,18.926268798813,170.336419189317,0.5169237685224889,0.5169237685224889),
    @("sCs","Bmp6","Bmpr1b","sCs",3,1,14.640471,43.921413,0.5728138386745201,0.5728138386745202,2,0.6666666666666666,0.1397713333333333,0.419314,0.09757108920650338,0.09757108920650338,2.046318152298,18.416863370682,0.05589007015203123,0.05589007015203125)
)

$rowIndex = 2
foreach ($row in $data) {
    $colIndex = 1
    foreach ($val in $row) {
        $ws.Cells.Item($rowIndex, $colIndex).Value = $val
        $colIndex++
    }
    $rowIndex++
}
